$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Agosto de 2020 a las 00:30'

$data = @(
  @('Estados Unidos', 5994417, 38689, 3281238, 2529690, 0, 1125, 183489),
  @('Brasil', 3717156, 42980, 2908848, 690643, 0, 999, 117665),
  @('India', 3307749, 75995, 2523443, 723677, 0, 1017, 60629),
  @('Rusia', 970865, 4676, 786150, 168032, 0, 115, 16683),
  @('Sudafrica', 615701, 2684, 525242, 76957, 0, 194, 13502),
  @('Peru', 607382, 0, 414577, 164804, 0, 0, 28001),
  @('Mexico', 568621, 4916, 393101, 114070, 0, 650, 61450),
  @('Colombia', 562128, 0, 395470, 148769, 0, 0, 17889),
  @('España', 426818, 3594, 0, 0, 0, 47, 28971),
  @('Chile', 402365, 1380, 376268, 15107, 0, 32, 10990),
  @('Iran', 365606, 2243, 314870, 29716, 0, 119, 21020),
  @('Argentina', 359638, 0, 268801, 83176, 0, 98, 7661),
  @('Reino Unido', 328846, 1048, 0, 0, 0, 16, 41465),
  @('Arabia Saudita', 310836, 1068, 284945, 22136, 0, 33, 3755),
  @('Banglades', 302147, 2519, 190183, 107882, 0, 54, 4082),
  @('Pakistan', 294193, 482, 278939, 8987, 0, 12, 6267),
  @('Italia', 262540, 1367, 206329, 20753, 0, 13, 35458),
  @('Turquia', 262507, 1313, 239797, 16527, 0, 20, 6183),
  @('Francia', 253587, 5429, 85524, 137519, 0, 0, 30544),
  @('Alemania', 239000, 1428, 209600, 20048, 0, 7, 9352),
  @('Irak', 215784, 3837, 157215, 51901, 0, 72, 6668),
  @('Filipinas', 202361, 5277, 133460, 65764, 0, 99, 3137),
  @('Indonesia', 160165, 2306, 115409, 37812, 0, 86, 6944),
  @('Canada', 126417, 448, 112455, 4868, 0, 4, 9094),
  @('Catar', 117742, 244, 114558, 2990, 0, 0, 194),
  @('Bolivia', 110999, 851, 48875, 57460, 0, 86, 4664),
  @('Ecuador', 110549, 1519, 95097, 9042, 0, 42, 6410),
  @('Ucrania', 110085, 1670, 53454, 54277, 0, 36, 2354),
  @('Israel', 108403, 1943, 86466, 21062, 0, 16, 875),
  @('Kazajistan', 105075, 173, 93990, 9562, 0, 0, 1523),
  @('Egipto', 97825, 206, 68713, 23795, 0, 19, 5317),
  @('Republica Dominicana', 92557, 340, 63478, 27466, 0, 28, 1613),
  @('Panama', 88381, 0, 62759, 23703, 0, 0, 1919),
  @('Suecia', 87072, 0, 0, 0, 0, 5, 5817),
  @('China', 84996, 15, 80015, 347, 0, 0, 4634),
  @('Oman', 84818, 166, 79409, 4763, 0, 4, 646),
  @('Belgica', 82447, 355, 18291, 54278, 0, 0, 9878),
  @('Kuwait', 82271, 698, 73906, 7844, 0, 2, 521),
  @('Rumania', 81646, 1256, 36286, 41939, 0, 54, 3421),
  @('Bielorrusia', 70974, 247, 69378, 939, 0, 5, 657),
  @('Guatemala', 70714, 1063, 58783, 9269, 0, 32, 2662),
  @('Paises Bajos', 68114, 571, 0, 0, 0, 8, 6215),
  @('Emiratos Arabes Unidos', 68020, 399, 59070, 8572, 0, 1, 378),
  @('Japon', 63822, 701, 51688, 10925, 0, 13, 1209),
  @('Polonia', 63802, 729, 43399, 18409, 0, 17, 1994),
  @('Singapur', 56495, 60, 54971, 1497, 0, 0, 27),
  @('Portugal', 56274, 362, 41184, 13283, 0, 2, 1807),
  @('Honduras', 55877, 398, 9124, 45050, 0, 20, 1703),
  @('Marruecos', 55864, 1336, 40586, 14294, 0, 29, 984),
  @('Nigeria', 52800, 0, 39964, 11829, 0, 0, 1007),
  @('Barein', 50393, 317, 47049, 3158, 0, 0, 186),
  @('Etiopia', 45221, 1533, 16311, 28185, 0, 16, 725),
  @('Ghana', 43769, 52, 42048, 1451, 0, 0, 270),
  @('Kirguistan', 43358, 113, 37217, 5084, 0, 0, 1057),
  @('Armenia', 43067, 131, 36726, 5480, 0, 3, 861),
  @('Argelia', 42619, 391, 29886, 11268, 0, 9, 1465),
  @('Venezuela', 41158, 0, 32015, 8800, 0, 0, 343),
  @('Suiza', 40645, 383, 34800, 3842, 0, 1, 2003),
  @('Uzbekistan', 39964, 323, 36402, 3271, 0, 7, 291),
  @('Afganistan', 38113, 43, 29042, 7670, 0, 4, 1401),
  @('Costa Rica', 36307, 1002, 13317, 22604, 0, 10, 386),
  @('Azerbaiyan', 35707, 148, 33281, 1904, 0, 1, 522),
  @('Moldavia', 34982, 624, 24156, 9859, 0, 7, 967),
  @('Nepal', 34418, 885, 19504, 14739, 0, 11, 175),
  @('Kenia', 33016, 213, 19296, 13156, 0, 5, 564),
  @('Serbia', 30974, 154, 29343, 924, 0, 2, 707),
  @('Irlanda', 28363, 162, 23364, 3222, 0, 0, 1777),
  @('Austria', 26033, 327, 22145, 3155, 0, 0, 733),
  @('Australia', 25204, 151, 20100, 4555, 0, 24, 549),
  @('El Salvador', 25140, 154, 13046, 11407, 0, 9, 687),
  @('Chequia', 22790, 242, 16939, 5433, 0, 2, 418),
  @('Estado de Palestina', 20155, 477, 13929, 6089, 0, 4, 137),
  @('Camerun', 18973, 0, 16540, 2023, 0, 0, 410),
  @('Bosnia y Herzegovina', 18609, 283, 12336, 5702, 0, 11, 571),
  @('Corea del Sur', 18265, 320, 14368, 3585, 0, 2, 312),
  @('Costa de Marfil', 17603, 41, 15941, 1548, 0, 0, 114),
  @('Dinamarca', 16537, 57, 14603, 1311, 0, 0, 623),
  @('Bulgaria', 15751, 162, 10750, 4415, 0, 14, 586),
  @('Madagascar', 14554, 79, 13582, 791, 0, 3, 181),
  @('Libano', 14248, 561, 3955, 10154, 0, 1, 139),
  @('Paraguay', 14228, 0, 7883, 6114, 0, 0, 231),
  @('Republica de Macedonia', 13914, 115, 10507, 2829, 0, 5, 578),
  @('Senegal', 13186, 130, 8852, 4059, 0, 1, 275),
  @('Sudan', 12974, 0, 6557, 5598, 0, 0, 819),
  @('Libia', 11834, 553, 1152, 10472, 0, 7, 210),
  @('Zambia', 11376, 91, 10693, 401, 0, 0, 282),
  @('Noruega', 10504, 50, 9150, 1090, 0, 0, 264),
  @('Consejo Danes para los Refugiados', 9912, 21, 8987, 671, 0, 3, 254),
  @('Malasia', 9291, 6, 8978, 188, 0, 0, 125),
  @('Grecia', 9280, 293, 3804, 5228, 0, 5, 248),
  @('Guinea', 9167, 39, 8150, 960, 0, 0, 57),
  @('Guayana Francesa', 8936, 32, 8461, 419, 0, 0, 56),
  @('Albania', 8927, 168, 4633, 4031, 0, 4, 263),
  @('Croacia', 8888, 358, 6362, 2351, 0, 0, 175),
  @('Gabon', 8468, 59, 7066, 1349, 0, 0, 53),
  @('Tayikistan', 8413, 34, 7214, 1132, 0, 0, 67),
  @('Haiti', 8122, 10, 5677, 2248, 0, 1, 197),
  @('Finlandia', 8002, 21, 7200, 467, 0, 0, 335),
  @('Luxemburgo', 7928, 90, 7140, 664, 0, 0, 124),
  @('Maldivas', 7225, 178, 4561, 2636, 0, 0, 28),
  @('Mauritania', 6977, 17, 6356, 463, 0, 0, 158),
  @('Namibia', 6431, 271, 2734, 3638, 0, 2, 59),
  @('Zimbabue', 6251, 55, 5001, 1071, 0, 13, 179),
  @('Malaui', 5474, 51, 3085, 2216, 0, 3, 173),
  @('Republica de Yibuti', 5383, 0, 5297, 26, 0, 0, 60),
  @('Hungria', 5288, 73, 3734, 940, 0, 0, 614),
  @('Guinea Ecuatorial', 4928, 2, 3798, 1047, 0, 0, 83),
  @('Hong Kong', 4736, 25, 4161, 497, 0, 1, 78),
  @('Republica de Africa Central', 4698, 7, 1778, 2859, 0, 0, 61),
  @('Nicaragua', 4494, 0, 2913, 1444, 0, 0, 137),
  @('Montenegro', 4444, 0, 3480, 877, 0, 0, 87),
  @('Suazilandia', 4327, 0, 2959, 1282, 0, 0, 86),
  @('Congo', 3979, 0, 1742, 2159, 0, 0, 78),
  @('Cuba', 3759, 15, 3159, 508, 0, 1, 92),
  @('Surinam', 3698, 0, 2818, 819, 0, 0, 61),
  @('Cabo Verde', 3630, 62, 2713, 879, 0, 1, 38),
  @('Ruanda', 3625, 88, 1810, 1800, 0, 0, 15),
  @('Mozambique', 3590, 82, 1927, 1642, 0, 0, 21),
  @('Eslovaquia', 3536, 84, 2192, 1311, 0, 0, 33),
  @('Tailandia', 3403, 1, 3237, 108, 0, 0, 58),
  @('Somalia', 3275, 0, 2443, 737, 0, 0, 95),
  @('Mayotte', 3237, 0, 2964, 234, 0, 0, 39),
  @('Tunez', 3069, 0, 1456, 1542, 0, 0, 71),
  @('Sri Lanka', 2984, 13, 2819, 153, 0, 0, 12),
  @('Lituania', 2726, 32, 1794, 847, 0, 0, 85),
  @('Eslovenia', 2722, 36, 2170, 419, 0, 0, 133),
  @('Mali', 2717, 4, 2052, 539, 0, 1, 126),
  @('Gambia', 2708, 22, 611, 2004, 0, 3, 93),
  @('Uganda', 2524, 98, 1268, 1230, 0, 1, 26),
  @('Sudan del Sur', 2510, 3, 1290, 1173, 0, 0, 47),
  @('Siria', 2440, 75, 550, 1792, 0, 3, 98),
  @('Angola', 2332, 49, 977, 1252, 0, 1, 103),
  @('Estonia', 2311, 17, 2054, 193, 0, 0, 64),
  @('Guinea-Bisau', 2205, 56, 1127, 1044, 0, 1, 34),
  @('Benin', 2145, 30, 1738, 367, 0, 1, 40),
  @('Islandia', 2082, 5, 1957, 115, 0, 0, 10),
  @('Sierra Leona', 2003, 2, 1577, 357, 0, 0, 69),
  @('Yemen', 1930, 6, 1097, 273, 0, 3, 560),
  @('Bahamas', 1813, 15, 634, 1131, 0, 2, 48),
  @('Jordania', 1756, 40, 1355, 386, 0, 1, 15),
  @('Malta', 1751, 46, 1077, 664, 0, 0, 10),
  @('Jamaica', 1732, 120, 840, 873, 0, 3, 19),
  @('Nueva Zelanda', 1695, 5, 1539, 134, 0, 0, 22),
  @('Aruba', 1670, 0, 565, 1097, 0, 0, 8),
  @('Botsuana', 1633, 71, 224, 1403, 0, 3, 6),
  @('Uruguay', 1536, 0, 1309, 184, 0, 0, 43),
  @('Republica de Chipre', 1484, 10, 935, 529, 0, 0, 20),
  @('Georgia', 1436, 7, 1150, 267, 0, 0, 19),
  @('Trinidad yTobago', 1384, 132, 192, 1177, 0, 0, 15),
  @('Reunion', 1372, 80, 692, 674, 0, 0, 6),
  @('Letonia', 1360, 18, 1135, 192, 0, 0, 33),
  @('Burkina Faso', 1352, 14, 1058, 239, 0, 0, 55),
  @('Togo', 1326, 17, 946, 353, 0, 0, 27),
  @('Liberia', 1298, 3, 832, 384, 0, 0, 82),
  @('Niger', 1173, 0, 1084, 20, 0, 0, 69),
  @('Principado de Andorra', 1098, 38, 893, 152, 0, 0, 53),
  @('Guyana', 1093, 33, 534, 528, 0, 0, 31),
  @('Lesoto', 1051, 2, 526, 494, 0, 1, 31),
  @('Vietnam', 1034, 5, 632, 373, 0, 2, 29),
  @('Republica del Chad', 998, 3, 873, 48, 0, 0, 77),
  @('Guadalupe', 935, 0, 289, 631, 0, 0, 15),
  @('Santo Tome y Principe', 892, 0, 833, 44, 0, 0, 15),
  @('Belice', 760, 30, 54, 695, 0, 1, 11),
  @('Crucero', 712, 0, 651, 48, 0, 0, 13),
  @('San Marino', 710, 0, 658, 10, 0, 0, 42),
  @('Birmania', 580, 76, 345, 229, 0, 0, 6),
  @('Tanzania', 509, 0, 183, 305, 0, 0, 21),
  @('Taiwan', 487, 0, 462, 18, 0, 0, 7),
  @('Martinica', 464, 0, 98, 350, 0, 0, 16),
  @('Islas Turcas y Caicos', 431, 48, 102, 327, 0, 0, 2),
  @('Burundi', 430, 0, 345, 84, 0, 0, 1),
  @('Papua Nueva Guinea', 419, 0, 232, 183, 0, 0, 4),
  @('San Martin (Parte Holandesa)', 418, 10, 148, 253, 0, 0, 17),
  @('Comoras', 417, 0, 399, 11, 0, 0, 7),
  @('Islas Feroe', 411, 0, 357, 54, 0, 0, 0),
  @('Polinesia Francesa', 372, 0, 148, 224, 0, 0, 0),
  @('Mauricio', 348, 0, 335, 3, 0, 0, 10),
  @('Isla de Man', 336, 0, 312, 0, 0, 0, 24),
  @('Eritrea', 315, 0, 276, 39, 0, 0, 0),
  @('Mongolia', 300, 2, 289, 11, 0, 0, 0),
  @('Camboya', 273, 0, 264, 9, 0, 0, 0),
  @('Gibraltar', 270, 14, 203, 67, 0, 0, 0),
  @('Islas Caimanes', 205, 0, 202, 2, 0, 0, 1),
  @('San Martin (Parte Francesa)', 198, 0, 52, 141, 0, 0, 5),
  @('Butan', 173, 17, 118, 55, 0, 0, 0),
  @('Bermudas', 168, 0, 151, 8, 0, 0, 9),
  @('Barbados', 164, 0, 132, 25, 0, 0, 7),
  @('Brunei', 144, 0, 139, 2, 0, 0, 3),
  @('Seychelles', 136, 0, 127, 9, 0, 0, 0),
  @('Monaco', 121, 0, 85, 35, 0, 0, 1),
  @('Liechtenstein', 102, 0, 94, 7, 0, 0, 1),
  @('Antigua y Barbuda', 94, 0, 89, 2, 0, 0, 3),
  @('San Vicente y las Granadinas', 58, 0, 57, 1, 0, 0, 0),
  @('Curazao', 49, 2, 34, 14, 0, 0, 1),
  @('Macao', 46, 0, 46, 0, 0, 0, 0),
  @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
  @('Guam', 32, 0, 0, 31, 0, 0, 1),
  @('Fiyi', 28, 0, 23, 4, 0, 0, 1),
  @('Islas Virgenes Britanicas', 26, 0, 8, 17, 0, 0, 1),
  @('Santa Lucia', 26, 0, 25, 1, 0, 0, 0),
  @('Timor Oriental', 26, 0, 25, 1, 0, 0, 0),
  @('Granada', 24, 0, 24, 0, 0, 0, 0),
  @('Nueva Caledonia', 23, 0, 23, 0, 0, 0, 0),
  @('Laos', 22, 0, 21, 1, 0, 0, 0),
  @('Dominica', 20, 0, 18, 2, 0, 0, 0),
  @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
  @('San Bartolome', 17, 0, 9, 8, 0, 0, 0),
  @('San Cristobal y Nieves', 17, 0, 17, 0, 0, 0, 0),
  @('Groenlandia', 14, 0, 14, 0, 0, 0, 0),
  @('Bonaire, San Eustaquio y Saba', 13, 0, 7, 6, 0, 0, 0),
  @('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0),
  @('Montserrat', 13, 0, 12, 0, 0, 0, 1),
  @('Santa Sede', 12, 0, 12, 0, 0, 0, 0),
  @('Sahara Occidental', 10, 0, 8, 1, 0, 0, 1),
  @('San Pedro y Miquelon', 5, 0, 1, 4, 0, 0, 0),
  @('Anguila', 3, 0, 3, 0, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
